$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Headers
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy formatting: bold header style from an existing sheet header cell,
# and the date number format from an existing date column cell.
$wsWeekly.Range("A1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A15").PasteSpecial(-4122)

# Data rows
$newSheet.Range("A2").Value = 45550.99999999999
$newSheet.Range("B2").Value = 601
$newSheet.Range("C2").Value = 297.2113874824213
$newSheet.Range("D2").Value = 881.9876868097449

$newSheet.Range("A3").Value = 45578.99999999999
$newSheet.Range("B3").Value = 383
$newSheet.Range("C3").Value = 91.58111398016173
$newSheet.Range("D3").Value = 675.8082978435361

$newSheet.Range("A4").Value = 45585.99999999999
$newSheet.Range("B4").Value = 329
$newSheet.Range("C4").Value = 39.38237388400116
$newSheet.Range("D4").Value = 612.3676855113272

$newSheet.Range("A5").Value = 45592.99999999999
$newSheet.Range("B5").Value = 274
$newSheet.Range("C5").Value = 1.873958786468395
$newSheet.Range("D5").Value = 578.6531050665579

$newSheet.Range("A6").Value = 45599.99999999999
$newSheet.Range("B6").Value = 220
$newSheet.Range("C6").Value = -72.85963747144686
$newSheet.Range("D6").Value = 506.3097789380308

$newSheet.Range("A7").Value = 45606.99999999999
$newSheet.Range("B7").Value = 165
$newSheet.Range("C7").Value = -105.5467925562958
$newSheet.Range("D7").Value = 464.8153920532259

$newSheet.Range("A8").Value = 45613.99999999999
$newSheet.Range("B8").Value = 111
$newSheet.Range("C8").Value = -189.4356322248027
$newSheet.Range("D8").Value = 387.6682313711508

$newSheet.Range("A9").Value = 45620.99999999999
$newSheet.Range("B9").Value = 56
$newSheet.Range("C9").Value = -242.7907371781852
$newSheet.Range("D9").Value = 336.5384731613977

$newSheet.Range("A10").Value = 45627.99999999999
$newSheet.Range("B10").Value = 2
$newSheet.Range("C10").Value = -282.3569586720301
$newSheet.Range("D10").Value = 290.5581025359397

$newSheet.Range("A11").Value = 45634.99999999999
$newSheet.Range("B11").Value = 0
$newSheet.Range("C11").Value = -348.4947370126546
$newSheet.Range("D11").Value = 265.9435826666729

$newSheet.Range("A12").Value = 45641.99999999999
$newSheet.Range("B12").Value = 0
$newSheet.Range("C12").Value = -407.0879350290161
$newSheet.Range("D12").Value = 186.3106609885777

$newSheet.Range("A13").Value = 45648.99999999999
$newSheet.Range("B13").Value = 0
$newSheet.Range("C13").Value = -459.6511716295239
$newSheet.Range("D13").Value = 153.6382634498628

$newSheet.Range("A14").Value = 45655.99999999999
$newSheet.Range("B14").Value = 0
$newSheet.Range("C14").Value = -494.1030109113238
$newSheet.Range("D14").Value = 57.16155397264841

$newSheet.Range("A15").Value = 45662.99999999999
$newSheet.Range("B15").Value = 0
$newSheet.Range("C15").Value = -558.8476629532458
$newSheet.Range("D15").Value = 5.503801498705254
